# Update base model version to match final EPS 2.1.1-us-v2
# (CPI input data refreshed from 2018 release to 2019 release)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# --- "About" sheet -------------------------------------------------------
# Year label 2018 -> 2019
$wsAbout.Range("B4").Value = 2019

# Source link text + target: 2018 PDF -> 2019 PDF
$newUrl = "https://www.bls.gov/cpi/tables/supplemental-files/historical-cpi-u-201912.pdf"
$wsAbout.Range("B6").Value = $newUrl
$wsAbout.Hyperlinks.Item(1).Address = $newUrl

# Make "About" the active/selected sheet (was "Data")
$wsAbout.Activate()

# --- "Data" sheet ----------------------------------------------------------
# Append the new 2019 annual-average CPI-U row (row 57), following the same
# pattern as the existing historical rows.
$wsData.Range("A57").Value = "2019.............................................................................     ."
$wsData.Range("B57").Value = 254.412
$wsData.Range("C57").Value = 256.903
$wsData.Range("D57").Value = 255.657
$wsData.Range("E57").Value = 2.3
$wsData.Range("F57").Value = 1.8
$wsData.Range("G57").Formula = '=$D$50/D57'
$wsData.Range("G57").NumberFormat = "0.000"
